$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix missing 를 particle in D29 (뇌내성) and D65 (빙내성) effect strings
$ws.Range("D29").Value = '뇌내성+20. 합계내성이 25이상이면 뇌속성 피해 소와 대를 무효화/뇌내성+15. 합계내성이 15이상이면 뇌속성 피해 소를 무효화/뇌내성-20'
$ws.Range("D65").Value = '빙내성+20. 합계내성이 25이상이면 빙속성 피해 소와 대를 무효화/빙내성+15. 합계내성이 15이상이면 빙속성 피해 소를 무효화/빙내성-20'

# Translate remaining Japanese skill rows (171-206) to Korean (with Japanese in parentheses)
$ws.Range("A171").Value = '발도회심(抜刀会心)'
$ws.Range("B171").Value = '발도술【기】(抜刀術【技】)'
$ws.Range("D171").Value = '발도 공격 시에 회심률이 +100% 늘어난다.'
$ws.Range("A172").Value = '발도멸기(抜刀減気)'
$ws.Range("B172").Value = '발도술【력】(抜刀術【力】)'
$ws.Range("D172").Value = '절단속성의 발도 공격 시에 공격력5와 멸기치15(쌍검은5)가 추가된다. 머리 공격 시에는 기절치 30도 추가'
$ws.Range("A173").Value = '공복(腹減り)'
$ws.Range("B173").Value = '공복 무효(腹減り無効)/공복 반감(腹減り半減)/공복 배가【소】(腹減り倍加【小】)/공복 배가【대】(腹減り倍加【大】)'
$ws.Range("D173").Value = '공복에 의한 스태미너 감소가 없어진다./공복에 의한 스태미너 감소가 1/2이 된다./공복에 의한 스태미너 감소속도가 1.5배로 올라간다./공복에 의한 스태미너 감소속도가 2배로 올라간다.'
$ws.Range("A174").Value = '반동(反動)'
$ws.Range("B174").Value = '반동경감+3(反動軽減+3)/반동경감+2(反動軽減+2)/반동경감+1(反動軽減+1)/반동경감-1(反動軽減-1)/반동경감-2(反動軽減-2)/반동경감-3(反動軽減-3)'
$ws.Range("D174").Value = '보우건의 반동이 3단계 작아진다. 속사의 반동에는 효과가 없다./보우건의 반동이 2단계 작아진다. 속사의 반동에는 효과가 없다./보우건의 반동이 1단계 작아진다. 속사의 반동에는 효과가 없다./보우건의 반동이 1단계 커진다. 속사의 반동에는 효과가 없다./보우건의 반동이 2단계 커진다. 속사의 반동에는 효과가 없다./보우건의 반동이 3단계 커진다. 속사의 반동에는 효과가 없다.'
$ws.Range("A175").Value = '비행술집(飛行酒場)'
$ws.Range("B175").Value = '비행술집의 마음(飛行酒場の心)'
$ws.Range("D175").Value = 'KO술(KO術) + 피리불기 명인(笛吹き名人)'
$ws.Range("A176").Value = '비술(秘術)'
$ws.Range("B176").Value = '스킬 가점+2(スキル加点+2)'
$ws.Range("D176").Value = '장비하고 있는 스킬 포인트의 합계치에 +2가 된다.'
$ws.Range("A177").Value = '화속성 공격(火属性攻撃)'
$ws.Range("B177").Value = '화속성 공격강화+2(火属性攻撃強化+2)/화속성 공격강화+1(火属性攻撃強化+1)/화속성 공격약화(火属性攻撃弱化)'
$ws.Range("D177").Value = '화속성의 공격이 1.1배+6이 된다./화속성의 공격이 1.05배+4이 된다./화속성의 공격이 0.75배가 된다.'
$ws.Range("A178").Value = '화내성(火耐性)'
$ws.Range("B178").Value = '화내성【대】(火耐性【大】)/화내성【소】(火耐性【小】)/화내성 약화(火耐性弱化)'
$ws.Range("D178").Value = '화내성+20. 합계내성이 25이상이면 화속성 피해 소와 대를 무효화/화내성+15. 합계내성이 15이상이면 화속성 피해 소를 무효화/화내성-20'
$ws.Range("A179").Value = '풍압(風圧)'
$ws.Range("B179").Value = '풍압【대】 무효(風圧【大】無効)/풍압【소】 무효(風圧【小】無効)'
$ws.Range("D179").Value = '풍압【대】와 풍압【소】의 영향을 받지 않게 된다. 용풍압을 풍압【대】의 효과로 경감/풍압【소】의 영향을 받지 않게 된다.'
$ws.Range("A180").Value = '피리(笛)'
$ws.Range("B180").Value = '피리불기 명인(笛吹き名人)'
$ws.Range("D180").Value = '수렵적 연주의 효과시간이 늘어나고 회복효과가 있는 선율은 회복량이 많은 쪽의 효과가 나오기 쉬워진다. 또한 피리 아이템이 1/2의 확률로 부숴지지 않는다.'
$ws.Range("A181").Value = '홍두(紅兜)'
$ws.Range("B181").Value = '홍두의 혼(紅兜の魂)'
$ws.Range("D181").Value = '앙심(逆恨み) + 집중(集中)'
$ws.Range("A182").Value = '베르나(ベルナ)'
$ws.Range("B182").Value = '베르나의 마음(ベルナの心)'
$ws.Range("D182").Value = '화내성【소】(火耐性【小】) + 공복 반감(腹減り半減)'
$ws.Range("A183").Value = '변칙사격(変則射撃)'
$ws.Range("B183").Value = '특정사격강화(特定射撃強化)'
$ws.Range("D183").Value = '보우건의 무기내장탄과 활의 곡사, 강사의 위력이 1.2배로 상승'
$ws.Range("A184").Value = '방어(防御)'
$ws.Range("B184").Value = '방어력UP【대】(防御力UP【大】)/방어력UP【중】(防御力UP【中】)/방어력UP【소】(防御力UP【小】)/방어력DOWN【소】(防御力DOWN【小】)/방어력DOWN【중】(防御力DOWN【中】)/방어력DOWN【대】(防御力DOWN【大】)'
$ws.Range("D184").Value = '방어력 1.06배+25/방어력 1.03배+20/방어력+15/방어력-10. 단 1미만이 되지 않는다./방어력 0.95배-10. 단 1미만이 되지 않는다./방어력 0.9배-10. 단 1미만이 되지 않는다.'
$ws.Range("A185").Value = '포술(砲術)'
$ws.Range("B185").Value = '포술왕(砲術王)/포술사(砲術師)'
$ws.Range("D185").Value = '발리스타, 건랜스의 포격, 용격포의 위력이 1.2배가 된다. 챠지 액스의 유탄병 장착 속성해방 베기의 위력이 1.35배, 철갑유탄, 용격탄, 연폭유탄, 대포탄의 위력이 1.3배가 된다. 용격포의 냉각시간이 120초에서 90초로 단축된다. 또 건랜스의 히트 게이지 빨강의 변동이 1/2가 된다./발리스타, 건랜스의 포격, 용격포의 위력이 1.1배가 된다. 챠지 액스의 유탄병 장착 속성해방 베기의 위력이 1.3배, 철갑유탄, 용격탄, 연폭유탄, 대포탄의 위력이 1.15배가 된다. 또 건랜스의 히트 게이지 빨강의 변동이 2/3가 된다.'
$ws.Range("A186").Value = '포획(捕獲)'
$ws.Range("B186").Value = '포획 명인(捕獲名人)/포획 달인(捕獲達人)'
$ws.Range("D186").Value = '몬스터 포획 시의 포획보수칸이 2~3칸에서 3~4칸으로 늘어난다./몬스터 포획 시의 포획보수칸이 2~3칸에서 3칸으로 늘어난다.'
$ws.Range("A187").Value = '북진낫토류(北辰納豆流)'
$ws.Range("B187").Value = '끈적끈적 검법(ネバネバ剣法)'
$ws.Range("D187").Value = '불굴(不屈) + 러너(ランナー) + 스태미너 탈취(スタミナ奪取)'
$ws.Range("A188").Value = '모쇄(矛砕)'
$ws.Range("B188").Value = '모쇄의 혼(矛砕の魂)'
$ws.Range("D188").Value = '업물(業物) + 반동경감+2(反動軽減+2) + 정령의 가호(精霊の加護)'
$ws.Range("A189").Value = '폿케(ポッケ)'
$ws.Range("B189").Value = '폿케의 마음(ポッケの心)'
$ws.Range("D189").Value = '빙내성【소】(氷耐性【小】) + 만복(まんぷく)'
$ws.Range("A190").Value = '잠재력(本気)'
$ws.Range("B190").Value = '힘의 해방+2(力の解放+2)/힘의 해방+1(力の解放+1)'
$ws.Range("D190").Value = '받은 합계 데미지가 180간격으로, 혹은 대형 몬스터에게 발견된 상태가 합계 5분간이 될 때 마다 발동. 1분30초간, 회심률이 50% 오르고, 스태미너 소비가 1/4이 된다. 스킬 도전자(挑戦者), 풀차지(フルチャージ)와는 중복되지 않는다./받은 합계 데미지가 180간격으로, 혹은 대형 몬스터에게 발견된 상태가 합계 5분간이 될 때 마다 발동. 1분30초간, 회심률이 30% 오르고, 스태미너 소비가 1/2이 된다. 스킬 도전자(挑戦者), 풀차지(フルチャージ)와는 중복되지 않는다.'
$ws.Range("A191").Value = '마비(麻痺)'
$ws.Range("B191").Value = '마비 무효(麻痺無効)/마비 배가(麻痺倍加)'
$ws.Range("D191").Value = '마비 상태가 되지 않는다./마비상태의 시간이 2배가 된다.'
$ws.Range("A192").Value = '마비병추가(麻痺瓶追加)'
$ws.Range("B192").Value = '마비병 추가(麻痺ビン追加)'
$ws.Range("D192").Value = '마비병의 장착이 가능하게 된다.'
$ws.Range("A193").Value = '수속성 공격(水属性攻撃)'
$ws.Range("B193").Value = '수속성 공격강화+2(水属性攻撃強化+2)/수속성 공격강화+1(水属性攻撃強化+1)/수속성 공격약화(水属性攻撃弱化)'
$ws.Range("D193").Value = '수속성의 공격이 1.1배+6이 된다./수속성의 공격이 1.05배+4가 된다./수속성의 공격이 0.75배가 된다.'
$ws.Range("A194").Value = '수내성(水耐性)'
$ws.Range("B194").Value = '수내성【대】(水耐性【大】)/수내성【소】(水耐性【小】)/수내성 약화(水耐性弱化)'
$ws.Range("D194").Value = '수내성+20. 합계내성이 25이상이면 수속성 피해 소와 대를 무효화/수내성+15. 합계내성이 15이상이면 수속성 피해 소를 무효화/수내성-20'
$ws.Range("A195").Value = '무상(無傷)'
$ws.Range("B195").Value = '풀 챠지(フルチャージ)'
$ws.Range("D195").Value = '체력이 최대일 때 공격력+20. 스킬 힘의해방(力の解放), 도전자(挑戦者)와 중복되지 않는다.'
$ws.Range("A196").Value = '무심(無心)'
$ws.Range("B196").Value = '명경지수(明鏡止水)'
$ws.Range("D196").Value = '수기 게이지의 축적량이 1.15배로 상승'
$ws.Range("A197").Value = '야초지식(野草知識)'
$ws.Range("B197").Value = '약초 초강화(薬草超強化)/약초 강화(薬草強化)'
$ws.Range("D197").Value = '약초의 회복력이 50이 된다. 광역화에 관해서는 스킬 소지자로부터 동료에 대한 변화는 없지만 동료로부터 스킬 소지자에 대한 회복량이 늘어난다./약초의 회복력이 30이 된다. 광역화에 관해서는 스킬 소지자로부터 동료에 대한 변화는 없지만 동료로부터 스킬 소지자에 대한 회복량이 늘어난다.'
$ws.Range("A198").Value = '유쿠모(ユクモ)'
$ws.Range("B198").Value = '유쿠모의 마음(ユクモの心)'
$ws.Range("D198").Value = '수내성【소】(水耐性【小】) + 허니 헌터(ハニーハンター)'
$ws.Range("A199").Value = '개열(鎧裂)'
$ws.Range("B199").Value = '개열의 혼(鎧裂の魂)'
$ws.Range("D199").Value = '가드 성능+2(ガード性能+2) + 납도술(納刀術)'
$ws.Range("A200").Value = '용기(龍気)'
$ws.Range("B200").Value = '용기활성(龍気活性)'
$ws.Range("D200").Value = '체력이 2/3이하가 되면 용속성 피해상태가 되지만 모든 속성내성치가 50이 되고 공격력이 1.1배가 된다. 용속성 피해는 공격력이 돌아오면 회복된다.'
$ws.Range("A201").Value = '용식선(龍識船)'
$ws.Range("B201").Value = '용식선의 마음(龍識船の心)'
$ws.Range("D201").Value = '용내성【소】(龍耐性【小】) + 세균연구가(細菌研究家)'
$ws.Range("A202").Value = '용속성 공격(龍属性攻撃)'
$ws.Range("B202").Value = '용속성 공격강화+2(龍属性攻撃強化+2)/용속성 공격강화+1(龍属性攻撃強化+1)/용속성 공격약화(龍属性攻撃弱化)'
$ws.Range("D202").Value = '용속성의 공격이 1.1배+6이 된다./용속성의 공격이 1.05배+4가 된다./용속성의 공격이 0.75배가 된다.'
$ws.Range("A203").Value = '용내성(龍耐性)'
$ws.Range("B203").Value = '용내성【대】(龍耐性【大】)/용내성【소】(龍耐性【小】)/용내성 약화(龍耐性弱化)'
$ws.Range("D203").Value = '용내성+20. 합계내성이 25이상이면 용속성 피해 소와 대를 무효화/용내성+15. 합계내성이 15이상이면 용속성 피해 소를 무효화/용내성-20'
$ws.Range("A204").Value = '유탄추가(榴弾追加)'
$ws.Range("B204").Value = '철갑유탄 전LV 추가(徹甲榴弾全LV追加)/철갑유탄 LV1 추가(徹甲榴弾LV1追加)'
$ws.Range("D204").Value = '전LV의 철갑유탄을 사용할 수 있게 된다./LV1철갑유탄을 사용할 수 있게 된다.'
$ws.Range("A205").Value = '열상(裂傷)'
$ws.Range("B205").Value = '열상 무효(裂傷無効)/열상 배가(裂傷倍加)'
$ws.Range("D205").Value = '열상 상태가 되지 않는다./열상상태가 2배가 된다.'
$ws.Range("A206").Value = '연격(連撃)'
$ws.Range("B206").Value = '연격의 심득(連撃の心得)'
$ws.Range("D206").Value = '공격 명중 시 5초간 회심률이 25% 상승. 도중 연속으로 공격이 5회 명중하면 회심률 상승이 30%가 된다.'

# Update the saved view state to match: scrolled down near the bottom of the table,
# with the last row (207, one past the data) selected as in the authored workbook.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 196
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("D207").Select()

